# "Generate Report for Handoff"
#
# This localization-status report has three sheets: "Overview", "zh-cn" and
# "de-de". The handoff process has now produced real handoff files for both
# target languages, so:
#   - the overall Status changes from "Handoff transform failed" to
#     "Ready for handoff" everywhere it is shown (Overview!B2/C2 and the
#     per-language sheets' B2),
#   - each language sheet gets a "Latest Handoff File" (hyperlinked) and a
#     "Latest Handoff Datetime" for row 2 (the source markdown file),
#   - each language sheet's Handoff Reason for row 2 becomes "Include"
#     (row 3, the .localization-config row, stays "Ignored").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$commit = "66af2ca68e3bbbd78f84de0bdb488713059e4662"
$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$commit"

$newStatus = "Ready for handoff"

# --- Status column updates (same text everywhere "Handoff transform failed" used to be) ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$zhcn.Range("B2").Value = $newStatus
$dede.Range("B2").Value = $newStatus

# NOTE: this runtime's PowerShell does not bind named (-Param value) function
# arguments correctly for COM objects (and even for plain strings), so the
# helper below takes plain positional parameters.
function Set-HandoffRow2 {
    param($ws, $xlfName, $handoffDatetime, $langCode)

    # Rebuild the row's hyperlinks so the new "Latest Handoff File" link
    # (column C) lands between the existing A2 and A3 links, keeping the
    # same display text for the two links that already existed.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/5758bb67-cef7-48ca-8c47-cf2be1947d43.md", "", "", "5758bb67-cef7-48ca-8c47-cf2be1947d43.md") | Out-Null

    $xlfUrl = "$baseUrl/e2e/Loc/$langCode/$xlfName"
    $ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl, "", "", $xlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config") | Out-Null

    # Latest Handoff Datetime
    $ws.Range("D2").Value = $handoffDatetime

    # Handoff Reason: row2 (real source file) is now included in the handoff,
    # row3 (.localization-config) stays ignored.
    $ws.Range("H2").Value = "Include"
    $ws.Range("H3").Value = "Ignored"
}

Set-HandoffRow2 $zhcn "5758bb67-cef7-48ca-8c47-cf2be1947d43.bac3cf26d6de083baf775e951c82bd585937bd17.zh-cn.xlf" "2016-02-18 10:15:37" "zh-cn"
Set-HandoffRow2 $dede "5758bb67-cef7-48ca-8c47-cf2be1947d43.bac3cf26d6de083baf775e951c82bd585937bd17.de-de.xlf" "2016-02-18 10:15:48" "de-de"
